$d = $word.ActiveDocument

# Update the date heading
$d.Content.Find.Execute("2025-04-30 Wednesday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-05-01 Thursday", 2) | Out-Null

# Update the division-problem answer table (positional access, since
# several answer strings repeat elsewhere in the table and a plain
# Find/Replace would not be able to disambiguate them).
$t = $d.Tables(1)

$t.Cell(1, 1).Range.Text = "75÷6=12, 3"
$t.Cell(1, 2).Range.Text = "77÷2=38, 1"
$t.Cell(1, 3).Range.Text = "51÷3=17, 0"
$t.Cell(1, 4).Range.Text = "26÷2=13, 0"
$t.Cell(1, 5).Range.Text = "25÷6=4, 1"

$t.Cell(5, 1).Range.Text = "14÷7=2, 0"
$t.Cell(5, 2).Range.Text = "93÷6=15, 3"
$t.Cell(5, 3).Range.Text = "98÷9=10, 8"
$t.Cell(5, 4).Range.Text = "20÷6=3, 2"
$t.Cell(5, 5).Range.Text = "55÷8=6, 7"

$t.Cell(9, 1).Range.Text = "20÷8=2, 4"
$t.Cell(9, 2).Range.Text = "91÷7=13, 0"
$t.Cell(9, 3).Range.Text = "47÷8=5, 7"
$t.Cell(9, 4).Range.Text = "40÷5=8, 0"
$t.Cell(9, 5).Range.Text = "55÷8=6, 7"

$t.Cell(13, 1).Range.Text = "31÷2=15, 1"
$t.Cell(13, 2).Range.Text = "70÷7=10, 0"
$t.Cell(13, 3).Range.Text = "71÷6=11, 5"
$t.Cell(13, 4).Range.Text = "20÷3=6, 2"
$t.Cell(13, 5).Range.Text = "40÷5=8, 0"

$t.Cell(17, 1).Range.Text = "69÷3=23, 0"
$t.Cell(17, 2).Range.Text = "98÷9=10, 8"
$t.Cell(17, 3).Range.Text = "68÷2=34, 0"
$t.Cell(17, 4).Range.Text = "74÷2=37, 0"
$t.Cell(17, 5).Range.Text = "91÷3=30, 1"
